$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-23: update date serial value 45186 -> 45188
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
